$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Feature"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Status"

# --- Column A (IDs), interleaved like the original authoring ---
$ws.Range("A2").Value = "CR1"
$ws.Range("A3").Value = "CR2"

# --- Column B (Feature names) ---
$ws.Range("B2").Value = "Registeration"
$ws.Range("B3").Value = "login"
$ws.Range("B4").Value = "Admin"
$ws.Range("B5").Value = "Rating"

# --- back to column A ---
$ws.Range("A4").Value = "CR3"
$ws.Range("A5").Value = "CR4"
$ws.Range("A6").Value = "CR5"

# --- Column C (Descriptions) ---
$ws.Range("C2").Value = "The system shall allow new users to create an account by filling in personal information including `n1-username,`n 2-email`n3-phone number`n4-password`nThe username must contain at least 3 letters and may include special characters . Each registration is validated to ensure that the user data is unique and complies with formatting rules."
$ws.Range("C3").Value = "Users shall be able to log in using their registered email and password`n The password must be at least  8 characters long and include at least one special character.`n Successful authentication will redirect users to their dashboard."
$ws.Range("C4").Value = "The system shall provide an admin interface with permissions to `nadd, delete, and update users and their Feedback. `nThis feature helps maintain content quality and user management across the platform."

$ws.Range("B6").Value = "Booking"

$ws.Range("C5").Value = "Users shall be able to rate travel destinations using a`n 5-star rating system. `nThese ratings will be saved in the user’s history and displayed visually for feedback and reference by others."
$ws.Range("C6").Value = "The platform shall enable users to search for and book flights by selecting travel dates`n,destinations, and airlines. Once a booking is made, data will be collected and sent to a third-party service to check for flight availability and confirm the reservation."

# --- Column D (Status) ---
$ws.Range("D2").Value = "Approved "
$ws.Range("D3").Value = "Approved "

# --- Alignment / styles ---
$ws.Range("A1:B6").HorizontalAlignment = -4108
$ws.Range("D1").HorizontalAlignment = -4108

$ws.Range("C2:C3").WrapText = $true

$ws.Range("C5").HorizontalAlignment = -4131
$ws.Range("C5").VerticalAlignment = -4108
$ws.Range("C5").WrapText = $true

$ws.Range("C4").HorizontalAlignment = -4131
$ws.Range("C4").VerticalAlignment = -4160
$ws.Range("C4").WrapText = $true
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows(2).RowHeight = 150
$ws.Rows(3).RowHeight = 90
$ws.Rows(4).RowHeight = 75
$ws.Rows(5).RowHeight = 90
$ws.Rows(6).RowHeight = 90

# --- Column widths ---
$ws.Columns(2).ColumnWidth = 12.8
$ws.Columns(3).ColumnWidth = 47.7

# --- Data validation on column D ---
$range = $ws.Range("D1:D1048576")
$range.Validation.Add(3, 1, 1, '"Approved , Not Approved"')

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("J2").Select() | Out-Null
